$wb = $excel.ActiveWorkbook

# --- Update the "Logs" sheet: append row 7 with the new mail entry ---
$logs = $wb.Worksheets("Logs")

$logs.Range("A7").Value = "Sollicitatie marketingfunctie"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D7").Value = "Sollicitatie / Vacature"
$logs.Range("F7").Value = "2025-06-19 21:16:11"
$logs.Range("G7").Value = "Nee"

# Extend the conditional formatting ranges on "Logs" to cover the new row
$dFormats = $logs.Range("D2:D6").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D7"))
}

$gFormats = $logs.Range("G2:G6").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G7"))
}

# --- Update the "Dashboard" sheet: append row 7 with the category count ---
$dash = $wb.Worksheets("Dashboard")

$dash.Range("A7").Value = "Sollicitatie / Vacature"
$dash.Range("B7").Value = 1

# --- Update the chart series references to include the new row ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$7"
$series.Values = "='Dashboard'!`$B`$2:`$B`$7"
